$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new values look numeric need to be forced to text
# so Excel stores them the same way as the original (inline string),
# matching the source data (which is formatted/text, not numeric).
$textForceCells = @("D5","D7","D8","D9","D10","D11","D12","D13","D14","D17", `
  "D21","D22","D23","D24","D25","D27","D28","D29","D30","D31","D32","D33", `
  "D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46", `
  "D47","D49","D50","D51")
foreach ($addr in $textForceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.236.14"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.905.78"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
$ws.Range("D5").Value = "307.32"

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.07%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5260"
$ws.Range("E7").Value = "  +0.74%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3813"
$ws.Range("E8").Value = "  +1.30%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.07281"

# Row 10 - Solana
$ws.Range("D10").Value = "21.78"
$ws.Range("E10").Value = "  +2.65%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "0.9030"
$ws.Range("E11").Value = "  -0.27%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.08179"
$ws.Range("E12").Value = "  -3.31%  "

# Row 13 - Litecoin
$ws.Range("D13").Value = "96.27"
$ws.Range("E13").Value = "  -0.87%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "5.362"
$ws.Range("E14").Value = "  +1.14%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "1.445.13"
$ws.Range("E15").Value = "  -24.27%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.16%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.000008661"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18 - Avalanche
$ws.Range("E18").Value = "  +1.47%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.09%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "27.275.84"
$ws.Range("E20").Value = "  +0.14%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "5.119"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22 - Cosmos
$ws.Range("D22").Value = "10.83"
$ws.Range("E22").Value = "  +1.70%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "6.513"
$ws.Range("E23").Value = "  +1.10%  "

# Row 24 - Monero
$ws.Range("D24").Value = "150.16"
$ws.Range("E24").Value = "  +2.06%  "

# Row 25 - LidoDAOToken
$ws.Range("D25").Value = "2.309"
$ws.Range("E25").Value = "  -0.85%  "

# Row 26 - EthereumClassic
$ws.Range("E26").Value = "  -0.10%  "

# Row 27 - Toncoin
$ws.Range("D27").Value = "1.743"
$ws.Range("E27").Value = "  -0.90%  "

# Row 28 - BitcoinCash
$ws.Range("D28").Value = "116.72"
$ws.Range("E28").Value = "  +1.29%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "4.845"
$ws.Range("E29").Value = "  +0.46%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "4.858"
$ws.Range("E30").Value = "  -1.11%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "0.09243"
$ws.Range("E31").Value = "  -0.61%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "0.8356"
$ws.Range("E32").Value = "  +4.94%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.05061"
$ws.Range("E33").Value = "  -0.15%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -0.89%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.989"
$ws.Range("E35").Value = "  +1.44%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "2.728"
$ws.Range("E36").Value = "  +5.68%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "3.338"

# Row 38 - TheSandbox
$ws.Range("D38").Value = "0.5822"
$ws.Range("E38").Value = "  +1.52%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.02006"
$ws.Range("E39").Value = "  -0.24%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "1.079"
$ws.Range("E40").Value = "  +0.27%  "

# Row 41 - Aptos
$ws.Range("D41").Value = "9.172"
$ws.Range("E41").Value = "  +1.17%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "6.609"
$ws.Range("E42").Value = "  -0.07%  "

# Row 43 - Quant
$ws.Range("D43").Value = "117.48"
$ws.Range("E43").Value = "  +1.11%  "

# Row 44 - Algorand
$ws.Range("D44").Value = "0.1521"
$ws.Range("E44").Value = "  +0.18%  "

# Row 45 - Decentraland
$ws.Range("D45").Value = "0.4933"
$ws.Range("E45").Value = "  +1.48%  "

# Row 46 - swapped from EnergySwap to PaxDollar
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47 - swapped from PaxDollar to EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.16"
$ws.Range("E47").Value = "  -0.21%  "

# Row 48 - NEARProtocol
$ws.Range("E48").Value = "  +0.97%  "

# Row 49 - Elrond
$ws.Range("D49").Value = "38.91"
$ws.Range("E49").Value = "  +3.26%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.06143"
$ws.Range("E50").Value = "  +3.02%  "

# Row 51 - Aave
$ws.Range("D51").Value = "64.46"
$ws.Range("E51").Value = "  +0.65%  "
